# Auto-generated edit script
# Applies the diff changes to "VENTAS POR GRUPO" (sheet1) and "VENTA MENSUAL" (sheet2)

$wb = $excel.ActiveWorkbook
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$pairsSheet1 = @(
    @("L5",0),
    @("D6",0),
    @("E6",0),
    @("L6",0),
    @("M6",0),
    @("I12",0),
    @("K14",0),
    @("D20",0),
    @("C22",0),
    @("D22",0),
    @("L22",0),
    @("C24",0),
    @("H24",0),
    @("I24",0),
    @("M26",0),
    @("D27",0),
    @("M27",0),
    @("K28",0),
    @("M28",0),
    @("Q28",0),
    @("E29",0),
    @("M29",0),
    @("L31",0),
    @("C36",0),
    @("D36",0),
    @("E36",0),
    @("G36",0),
    @("L36",0),
    @("M36",0),
    @("Q39",0),
    @("M41",0),
    @("D42",0),
    @("D43",0),
    @("M43",0),
    @("N43",0),
    @("D44",0),
    @("M44",0),
    @("M45",0),
    @("D50",0),
    @("M50",0),
    @("C55",'0 de 53'),
    @("D55",'0 de 53'),
    @("E55",'0 de 53'),
    @("G55",'0 de 53'),
    @("H55",'0 de 53'),
    @("I55",'0 de 53'),
    @("K55",'0 de 53'),
    @("L55",'0 de 53'),
    @("M55",'0 de 53'),
    @("N55",'0 de 53'),
    @("Q55",'0 de 53')
)

foreach ($p in $pairsSheet1) {
    $wsGrupo.Range($p[0]).Value = $p[1]
}

$pairsSheet2 = @(
    @("C1",'abril'),
    @("D1",'mayo'),
    @("E1",'junio'),
    @("F1",'julio'),
    @("C5",3471.96),
    @("D5",4158.27),
    @("E5",5087.23),
    @("F5",0),
    @("C6",710.14),
    @("D6",1528.39),
    @("E6",1516.28),
    @("F6",0),
    @("C11",262.99),
    @("D11",0),
    @("E12",135.9),
    @("F12",0),
    @("D14",2169.75),
    @("E14",456.84),
    @("F14",0),
    @("C19",6506.43),
    @("D19",0),
    @("E20",259.49),
    @("F20",0),
    @("C22",607.78),
    @("D22",4141.42),
    @("E22",2733.67),
    @("F22",0),
    @("C24",1567.1),
    @("D24",61.78),
    @("E24",3114.83),
    @("F24",0),
    @("C25",129.6),
    @("D25",0),
    @("E26",727.97),
    @("F26",0),
    @("C27",29332.26),
    @("D27",6249.76),
    @("E27",36680.13),
    @("F27",0),
    @("C28",7533.56),
    @("D28",7315.29),
    @("E28",3247.24),
    @("F28",0),
    @("C29",7135.59),
    @("D29",3563.29),
    @("E29",1079.23),
    @("F29",0),
    @("C31",4280.48),
    @("D31",0),
    @("E31",2568.3),
    @("F31",0),
    @("C36",722.5700000000001),
    @("D36",5704.92),
    @("E36",14177.18),
    @("F36",0),
    @("C39",2673.89),
    @("D39",2403.41),
    @("E39",2172.6),
    @("F39",0),
    @("C41",0),
    @("D41",2689.09),
    @("E41",-39.67),
    @("F41",0),
    @("C42",15577.98),
    @("D42",0),
    @("E42",86.5),
    @("F42",0),
    @("C43",3716.29),
    @("D43",7574.03),
    @("E43",6905.1),
    @("F43",0),
    @("C44",418.61),
    @("D44",731.63),
    @("E44",3152.12),
    @("F44",0),
    @("C45",1091.58),
    @("D45",722.54),
    @("E45",158.83),
    @("F45",0),
    @("C47",969.61),
    @("D47",798),
    @("E47",0),
    @("C49",0),
    @("C50",4141.1),
    @("D50",4953.13),
    @("E50",5333.85),
    @("F50",0),
    @("C51",3336.39),
    @("D51",-11.75),
    @("E51",0),
    @("D54",144),
    @("E54",0),
    @("C55",94185.91),
    @("D55",54896.95),
    @("E55",89553.62),
    @("F55",0)
)

foreach ($p in $pairsSheet2) {
    $wsMensual.Range($p[0]).Value = $p[1]
}

# Column width changes on "VENTA MENSUAL": D 14->13, E 13->14, F 14->11
$wsMensual.Columns.Item(4).ColumnWidth = 12.17
$wsMensual.Columns.Item(5).ColumnWidth = 13.17
$wsMensual.Columns.Item(6).ColumnWidth = 10.17

Write-Host "Edit complete"
